$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Content from 7.jpg" heading becomes a plain paragraph reading "材料2",
#    and the leading "材料2" that used to prefix the following paragraph is
#    removed (it now lives in its own paragraph above).
#    This must run before the table is touched, because locating paragraphs
#    by index after a table-row deletion is unreliable in this runtime, so
#    we do the paragraph-index based lookup first while indices are stable.
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Content from 7.jpg*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $headingPara = $d.Paragraphs.Item($targetIndex)
    # Remove the whole paragraph (text + its own paragraph mark/formatting,
    # i.e. the Heading1 style goes away with it) -- this merges the
    # following (unstyled) paragraph up into this slot.
    $fullRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
    $fullRange.Delete()

    # Insert a fresh, un-styled paragraph "材料2" in front of what remains.
    $mergedPara = $d.Paragraphs.Item($targetIndex)
    $insertPoint = $d.Range($mergedPara.Range.Start, $mergedPara.Range.Start)
    $insertPoint.InsertBefore("材料2" + [char]13)
}

# Strip the now-redundant "材料2" prefix from the paragraph that used to
# start with it.
$d.Content.Find.Execute("材料2朝廷在故都", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "朝廷在故都", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Table restructuring (the "材料1" population table).
#    The original table has a stray placeholder row (just "例" markers) and
#    the data rows are mis-aligned by one column; the fixed table drops that
#    placeholder row and realigns every data row, adding the missing
#    "37.1%" figure to the last row.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Remove the placeholder row (was row 3: [empty] | 例 | [empty] | 例 | [empty]).
$t.Rows.Item(3).Delete()

# Row 3 (was "例 | 西汉 | 2470685 | 19.8% | 9985785") -> realigned + extended.
$t.Rows.Item(3).Cells.Item(1).Range.Text = "西汉"
$t.Rows.Item(3).Cells.Item(2).Range.Text = "2470685"
$t.Rows.Item(3).Cells.Item(3).Range.Text = "19.8%"
$t.Rows.Item(3).Cells.Item(4).Range.Text = "9985785"
$t.Rows.Item(3).Cells.Item(5).Range.Text = "80.2%"

# Row 4 (was "80.2% | 唐代 | 3920415 | 43.2% | 5148529") -> realigned.
$t.Rows.Item(4).Cells.Item(1).Range.Text = "唐代"
$t.Rows.Item(4).Cells.Item(2).Range.Text = "3920415"
$t.Rows.Item(4).Cells.Item(3).Range.Text = "43.2%"
$t.Rows.Item(4).Cells.Item(4).Range.Text = "5148529"
$t.Rows.Item(4).Cells.Item(5).Range.Text = "56.8%"

# Row 5 (was "56.8% | 北宋 | 11224760 | 62.9% | 6624296") -> realigned +
# gains the new "37.1%" cell value.
$t.Rows.Item(5).Cells.Item(1).Range.Text = "北宋"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "11224760"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "62.9%"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "6624296"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "37.1%"

# ---------------------------------------------------------------------------
# 3) Paragraph numbering fix just above the table ("材料1" / "15.").
#    Do the "15." -> "南方" replacement FIRST: once "材料1" becomes
#    "15.材料1" it would itself contain a "15." substring and get caught by
#    a later search for "15.".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("15.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "南方", 2) | Out-Null
$d.Content.Find.Execute("材料1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "15.材料1", 2) | Out-Null
